$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve each target cell as Text (matches original shared-string cell type)
# while only rewriting the numeric-looking string content, per the sharedStrings diff.
$updates = @{
    "C2" = "-6.7811680951745945"
    "D2" = "-0.18797690491407337"
    "E2" = "53.098332452860106"
    "R2" = "-75.28828459824439"
    "S2" = "-41.64387311722784"
    "T2" = "41.267919307399694"
    "U2" = "48.53956602737416"
    "C3" = "-0.3602560427839556"
    "D3" = "1.0978467476340938"
    "E3" = "42.505066876072696"
    "R3" = "-46.26287772834596"
    "S3" = "41.518064868806654"
    "T3" = "42.62616006194195"
    "U3" = "-39.322371373538466"
    "C4" = "33.59474241236954"
    "D4" = "51.3791413042813"
    "E4" = "44.697679144728774"
    "R4" = "73.28907794449847"
    "S4" = "-41.66839090358291"
    "T4" = "59.085660814902106"
    "U4" = "43.67262179366049"
    "C5" = "22.054481794032124"
    "D5" = "46.78843228579012"
    "E5" = "48.1090131792516"
    "R5" = "55.65905418973773"
    "S5" = "-61.01799158518947"
    "T5" = "45.36209725704591"
    "U5" = "48.21476731453433"
    "C6" = "-7.206270743291331"
    "D6" = "-5.009749617147527"
    "E6" = "46.69603663169677"
    "R6" = "-63.8295716297941"
    "S6" = "45.023987890923834"
    "T6" = "-42.2805949123608"
    "U6" = "32.26109567806574"
    "C7" = "-35.65015784612282"
    "D7" = "-64.11421157479182"
    "E7" = "50.77323540196945"
    "R7" = "-78.50780065830939"
    "S7" = "35.671538694732746"
    "T7" = "-64.11421157479182"
    "C8" = "6.131246248070719"
    "D8" = "39.48392088041285"
    "E8" = "44.00518345892174"
    "R8" = "-50.0149666393685"
    "S8" = "-45.3720306789774"
    "T8" = "39.48392088041285"
    "U8" = "44.84279526031155"
    "V8" = "41.716512417975096"
    "C9" = "-13.341634463113877"
    "D9" = "-12.245968813072313"
    "E9" = "47.96643435536334"
    "R9" = "30.729645134302352"
    "S9" = "-71.04512852395577"
    "T9" = "-52.51097228210848"
    "U9" = "28.019034655963857"
    "V9" = "-58.72067016588287"
    "W9" = "43.47828440299765"
    "C10" = "-24.193722390095967"
    "D10" = "-47.604896876944316"
    "E10" = "51.45740053082853"
    "R10" = "-92.16758144775933"
    "S10" = "-53.03325700797723"
    "T10" = "38.378194693492325"
    "U10" = "-47.604896876944316"
    "V10" = "33.45892868870871"
    "C11" = "0.6338835804533058"
    "D11" = "1.100474317491166"
    "E11" = "39.32581387373229"
    "S11" = "37.9726530596958"
    "T11" = "41.79716546673732"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"